$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Internal Validity" sub-section was added under "Threats to Validity",
#     between "Construct Validity" and "Conclusion Validity" (new row 36;
#     everything from the old row 36 down shifts to row+1). ---
$ws.Rows(36).Insert()
$ws.Range("A36").Value = "Internal Validity"

# --- New column B: a page-count / weight figure next to every outline entry. ---
$ws.Range("B1").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("B3").Value = 4
$ws.Range("B4").Value = 4
$ws.Range("B5").Value = 1
$ws.Range("B6").Value = 1
$ws.Range("B7").Value = 1
$ws.Range("B8").Value = 1
$ws.Range("B9").Value = 1
$ws.Range("B10").Value = 0.5
$ws.Range("B11").Value = 3
$ws.Range("B12").Value = 3
$ws.Range("B13").Value = 0.5
$ws.Range("B14").Value = 0.5
$ws.Range("B15").Value = 0.5
$ws.Range("B16").Value = 1.33
$ws.Range("B17").Value = 1.33
$ws.Range("B18").Value = 1.34
$ws.Range("B19").Value = 0.5
$ws.Range("B20").Value = 2
$ws.Range("B21").Value = 3
$ws.Range("B22").Value = 0.5
$ws.Range("B23").Value = 3
$ws.Range("B24").Value = 0.5
$ws.Range("B25").Value = 0.5
$ws.Range("B26").Value = 2
$ws.Range("B27").Value = 3.5
$ws.Range("B28").Value = 0.5
# B29 ("Study Results") intentionally left blank, as in the source.
$ws.Range("B30").Value = 4.5
$ws.Range("B31").Value = 4
$ws.Range("B32").Value = 1.5
$ws.Range("B33").Value = 4
# B34 ("Threats to Validity") intentionally left blank, as in the source.
$ws.Range("B35").Value = 1
$ws.Range("B36").Value = 1
$ws.Range("B37").Value = 1
$ws.Range("B38").Value = 1
$ws.Range("B39").Value = 0.5
$ws.Range("B40").Value = 1
$ws.Range("B41").Value = 1
$ws.Range("B42").Value = 1
$ws.Range("B43").Value = 0.5
$ws.Range("B44").Value = 1
$ws.Range("B45").Value = 2
$ws.Range("B46").Value = 0.5
$ws.Range("B47").Value = 1.5
$ws.Range("B48").Value = 2
# B49 ("Study Results") intentionally left blank, as in the source.
$ws.Range("B50").Value = 1.5
$ws.Range("B51").Value = 1.5
$ws.Range("B52").Value = 1.5
# B53 ("CHAPTER 6 - CONCLUSION") intentionally left blank, as in the source.
$ws.Range("B54").Value = 1
$ws.Range("B55").Value = 1
$ws.Range("B56").Value = 1
$ws.Range("B57").Value = 1
$ws.Range("B58").Value = 1

# --- Two new summary rows appended at the bottom. ---
$ws.Range("A59").Value = "TOTAL"
$ws.Range("B59").Formula = "=SUM(B1:B58)"
$ws.Range("A60").Value = "TOTAL CONTENT"
$ws.Range("B60").Formula = "=SUM(B1:B4,B10:B56)"

# --- Column A was widened (best-fit style) to comfortably show the longest
#     heading now that the sheet carries real content next to it. ---
$ws.Columns("A").ColumnWidth = 88

# --- Selection/scroll position left where the author was last editing. ---
$ws.Range("D54").Select()
